$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.779.41"
$ws.Range("E2").Value = "  +4.99%  "
$ws.Range("D3").Value = "3.105.17"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'584.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("D6").Value = "'143.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.098.86"
$ws.Range("E8").Value = "  +2.78%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.17%  "
$ws.Range("D11").Value = "'5.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.80%  "
$ws.Range("D12").Value = "'0.467"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.51%  "
$ws.Range("D14").Value = "'35.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.70%  "
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "3.620.42"
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").Value = "'7.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "3.102.97"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").Value = "62.829.65"
$ws.Range("E19").Value = "  +5.09%  "
$ws.Range("D20").Value = "'463.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.78%  "
$ws.Range("D21").Value = "'14.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.77%  "
$ws.Range("D22").Value = "'0.729"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'7.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.08%  "
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "'81.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("E28").Value = "  +4.84%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  +5.05%  "
$ws.Range("D31").Value = "'6.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.64%  "
$ws.Range("D32").Value = "'27.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.64%  "
$ws.Range("E33").Value = "  +8.39%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0833"
$ws.Range("E34").Value = "  +5.31%  "
$ws.Range("B35").Value = "Stacks"
$ws.Range("C35").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D35").Value = "'2.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.69%  "
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("D38").Value = "'3.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.20%  "
$ws.Range("D39").Value = "'51.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.75%  "
$ws.Range("D40").Value = "'8.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("D41").Value = "'429.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.95%  "
$ws.Range("D42").Value = "2.914.40"
$ws.Range("E42").Value = "  +4.22%  "
$ws.Range("D43").Value = "'0.0369"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("D44").Value = "'0.278"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.49%  "
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("D46").Value = "'2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.94%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'34.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'123.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "'24.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.24%  "
